{"js": "// Load all paragraphs in the document body so we can locate the ones we\n// need to touch by their (current) text content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction indexOfText(text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + text);\n}\n\n// 1) The empty paragraph right after \"User Stories:\" becomes \"Cycle 2:\".\nconst titleIdx = indexOfText(\"User Stories:\");\nconst cycle2Para = items[titleIdx + 1];\ncycle2Para.insertText(\"Cycle 2:\", \"Replace\");\n\n// 2) The empty paragraph right after the \"...already dug holes.\" paragraph\n//    becomes \"Cycle 3:\" (this is where the new section header goes).\nconst enemiesIdx = indexOfText(\n  \"As a player, I would like to encounter at least two different types of enemies\\u2014one that can dig, and one that can only move through already dug holes.\"\n);\nconst cycle3Para = items[enemiesIdx + 1];\ncycle3Para.insertText(\"Cycle 3:\", \"Replace\");\nawait context.sync();\n\n// 3) Move the \"_GoBack\" bookmark out of the \"...lives that I can lose...\"\n//    paragraph (where it currently splits \"a\" / \"re publicly displayed.\")\n//    and place it between \"Cycle 3\" and \":\" instead.\nconst livesIdx = indexOfText(\n  \"As a player, I want to have lives that I can lose and gain and are publicly displayed.\"\n);\nconst livesPara = items[livesIdx];\nconst livesText = livesPara.text;\n\ncontext.document.deleteBookmark(\"_GoBack\");\n// Re-write the paragraph's text as a single run now that the bookmark\n// (which used to split it into two runs) is gone.\nlivesPara.insertText(livesText, \"Replace\");\nawait context.sync();\n\n// Split \"Cycle 3:\" into \"Cycle 3\" / \":\" with the bookmark re-inserted\n// right before the colon.\nconst cycle3Range = cycle3Para.getRange();\nconst colonResults = cycle3Range.search(\":\", { matchCase: true });\ncolonResults.load(\"items\");\nawait context.sync();\n\nconst colonStart = colonResults.items[0].getRange(\"Start\");\ncolonStart.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word is the Word Application; $d is the active document.\n$d = $word.ActiveDocument\n\n$wdCharacter = 1\n$wdCollapseStart = 1\n\nfunction Get-ParagraphIndexByText {\n    param(\n        [string]$ExactText\n    )\n    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        $t = $p.Range.Text\n        # Range.Text includes the trailing paragraph mark; strip it.\n        if ($t.Length -gt 0 -and [int][char]$t.Substring($t.Length - 1, 1) -eq 13) {\n            $t = $t.Substring(0, $t.Length - 1)\n        }\n        if ($t -eq $ExactText) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $ExactText\"\n}\n\n# 1) The empty paragraph right after \"User Stories:\" becomes \"Cycle 2:\".\n$titleIdx = Get-ParagraphIndexByText \"User Stories:\"\n$cycle2Idx = $titleIdx + 1\n$d.Paragraphs.Item($cycle2Idx).Range.Text = \"Cycle 2:\"\n\n# 2) The empty paragraph right after the \"...already dug holes.\" paragraph\n#    becomes \"Cycle 3:\" (new section header).\n$enemiesIdx = Get-ParagraphIndexByText \"As a player, I would like to encounter at least two different types of enemies$([char]0x2014)one that can dig, and one that can only move through already dug holes.\"\n$cycle3Idx = $enemiesIdx + 1\n$d.Paragraphs.Item($cycle3Idx).Range.Text = \"Cycle 3:\"\n\n# 3) Move the \"_GoBack\" bookmark out of the \"...lives that I can lose...\"\n#    paragraph (where it currently splits \"a\" / \"re publicly displayed.\")\n#    and place it between \"Cycle 3\" and \":\" instead.\n$livesIdx = Get-ParagraphIndexByText \"As a player, I want to have lives that I can lose and gain and are publicly displayed.\"\n\n# Delete the old bookmark.\n$d.Bookmarks.Item(\"_GoBack\").Delete()\n\n# Re-write the \"lives...\" paragraph's text as a single run (the bookmark\n# deletion leaves the old run split intact, so force a real text rewrite).\n$livesRange = $d.Paragraphs.Item($livesIdx).Range.Duplicate\n$livesRange.MoveEnd($wdCharacter, -1) | Out-Null\n$livesText = $livesRange.Text\n$livesRange.Text = \"\"\n$livesRange2 = $d.Paragraphs.Item($livesIdx).Range.Duplicate\n$livesRange2.MoveEnd($wdCharacter, -1) | Out-Null\n$livesRange2.Text = $livesText\n\n# Find the \":\" in the new \"Cycle 3:\" paragraph and add the bookmark right\n# before it, splitting \"Cycle 3\" and \":\" into separate runs.\n$cycle3Range = $d.Paragraphs.Item($cycle3Idx).Range.Duplicate\n$cycle3Range.Find.Execute(\":\") | Out-Null\n$cycle3Range.Collapse($wdCollapseStart)\n$d.Bookmarks.Add(\"_GoBack\", $cycle3Range) | Out-Null\n"}
